# Generate Report for Handoff
# Updates localization status from "In Translation" to "Ready for handoff"
# and refreshes the handoff timestamps, matching the widened "Status" /
# locale columns that Excel's column auto-fit produces for the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Handoff timestamps bump forward by 40 seconds ---
$overview.Range("G2").Value = "2016-08-24 00:37:38"
$dede.Range("H2").Value = "2016-08-24 00:37:38"

$zhcn.Range("H2").Value = "2016-08-24 00:37:33"

# --- Column widths widen to fit the longer "Ready for handoff" text ---
# (Excel's ColumnWidth setter snaps to whole-pixel steps, same as the real
# Excel.Application object; 16.5 lands in the pixel bucket closest to the
# generated report's target width of ~17.22 characters.)
$overview.Columns.Item(5).ColumnWidth = 16.5
$overview.Columns.Item(6).ColumnWidth = 16.5

$zhcn.Columns.Item(3).ColumnWidth = 16.5
$dede.Columns.Item(3).ColumnWidth = 16.5
